# Rename the "Index" header (column A, row 1) to "Id" and move the
# selection from G16 to A2, matching the committed workbook change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Id"

$ws.Range("A2").Select()
